$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting of the "purpose" column (E) values from fullRNASEQ to fullRNASeq
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Range("E$row")
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
